# ---------------------------------------------------------------------------
# Regenerate penyata to follow new data and format
# (3LUHUR-2023.xlsx : penyata-akhir-2023/form3)
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Relabel the four "Kali ..." checkpoints to "Semakan Kali ..." wherever
#    they occur (Merit Pendahuluan / Laporan Atas Talian / JPPM sections all
#    reuse the same four labels).
# ---------------------------------------------------------------------------
foreach ($r in 16,22,28) { $ws.Range("C$r").Value = "Semakan Kali Pertama" }
foreach ($r in 17,23,29) { $ws.Range("C$r").Value = "Semakan Kali Kedua" }
foreach ($r in 18,24,30) { $ws.Range("C$r").Value = "Semakan Kali Ketiga" }
foreach ($r in 19,25,31) { $ws.Range("C$r").Value = "Semakan Kali Keempat" }

# ---------------------------------------------------------------------------
# 2. Re-case the competition entry names (Penyertaan Pertandingan section).
# ---------------------------------------------------------------------------
$ws.Range("C34").Value = "Unga"
$ws.Range("C35").Value = "Bouquet Kreatif"
$ws.Range("C36").Value = "Kad Raya Untuk Guruku"
$ws.Range("C37").Value = "Riang Ria Kuih Raya"
$ws.Range("C38").Value = "Creative Collage"

# New competition entries added on the two previously-blank rows.
$ws.Range("C39").Value = "Teng Teng"
$ws.Range("C40").Value = "Doodle Koperasi"

# ---------------------------------------------------------------------------
# 3. Updated financial figures.
# ---------------------------------------------------------------------------
$ws.Range("D17").Value = 880
$ws.Range("E17").Value = 700
$ws.Range("D18").Value = 9670
$ws.Range("D39").Value = 100

# ---------------------------------------------------------------------------
# 4. Move the title "STATEMENT OF HOMEROOM ACCOUNT" from E4 to D4 and widen
#    its merge so it spans D4:G4.
# ---------------------------------------------------------------------------
$ws.Range("E4").Copy($ws.Range("D4"))
$ws.Range("E4").Clear()
$ws.Range("D4:G4").Merge()

# ---------------------------------------------------------------------------
# 5. New / widened merge regions that go with the refreshed layout.
# ---------------------------------------------------------------------------
$ws.Range("B5:C5").Merge()
$ws.Range("B12:F12").Merge()

$ws.Range("B21:C21").UnMerge()
$ws.Range("B21:E21").Merge()

$ws.Range("B27:C27").UnMerge()
$ws.Range("B27:E27").Merge()

$ws.Range("B33:C33").UnMerge()
$ws.Range("B33:E33").Merge()

$ws.Range("B15:C15").UnMerge()

$ws.Range("B43:E43").Merge()

# ---------------------------------------------------------------------------
# 6. Reposition / resize the logo picture (now anchored as a single-cell,
#    fixed-size picture instead of a two-cell stretch anchor).
# ---------------------------------------------------------------------------
$shp = $ws.Shapes.Item(1)
$shp.Placement = 3
$shp.Left = $ws.Columns.Item(1).Width + 12
$shp.Top = 14.25
$shp.Width = 46.5
$shp.Height = 47.25

# ---------------------------------------------------------------------------
# 7. Drop the trailing blank formatted row so the sheet ends at row 1000.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1001).Delete()

# ---------------------------------------------------------------------------
# 8. Sheet / workbook view clean-up matching the refreshed export.
# ---------------------------------------------------------------------------
$ws.Range("A1").Select()
$ws.PageSetup.CenterHorizontally = $true
$ws.PageSetup.FitToPagesTall = 1
$ws.PageSetup.FitToPagesWide = 1
